$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: matchsource -> matchsync
$wsMeta.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/ms-abo-group-codes"

# Title: MatchSource -> MatchSync
$wsMeta.Range("B5").Value = "MatchSync ABO Group Value Set"

# Experimental value cell (was blank) -> true
$wsMeta.Range("B7").Value = "true"

# Date bump
$wsMeta.Range("B8").Value = "2024-02-19T18:37:26-06:00"

# Description: MatchSource -> MatchSync
$wsMeta.Range("B11").Value = "MatchSync ABO group codes. Combines LOINC, Snomed, and NMDP codes"

# Include ValueSets sheet: nmdp-abo-codes URL matchsource -> matchsync
$wsInc1 = $wb.Worksheets.Item("Include ValueSets")
$wsInc1.Range("A2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-abo-codes"

# Include ValueSets 3 sheet: sct-abo-group-codes URL matchsource -> matchsync
$wsInc3 = $wb.Worksheets.Item("Include ValueSets 3")
$wsInc3.Range("A2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/sct-abo-group-codes"
